# Applies the data update described in the commit "Update figures with new data"
# to Sheet1 of the institution-stats workbook:
#  - Helsinki University Hospital (row 12): trials/percentage/CI figures revised
#  - Oerebro University / Oerebro University Hospital rows moved from the end of the
#    alphabetical institution list (after "Zealand University Hospital") to their
#    correct alphabetical position (right after "Odense University Hospital"),
#    shifting all institutions in between down by two rows
#  - Total row (row 56): trials and conf_int_ll updated to reflect the revised counts

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(12, 2).Value = 48
$ws.Cells.Item(12, 4).Value = 60.4
$ws.Cells.Item(12, 5).Value = 46.3
$ws.Cells.Item(12, 6).Value = 73
$ws.Cells.Item(27, 1).Value = "Örebro University"
$ws.Cells.Item(27, 2).Value = 17
$ws.Cells.Item(27, 3).Value = 6
$ws.Cells.Item(27, 4).Value = 35.3
$ws.Cells.Item(27, 5).Value = 17.3
$ws.Cells.Item(27, 6).Value = 58.7
$ws.Cells.Item(28, 1).Value = "Örebro University Hospital"
$ws.Cells.Item(28, 2).Value = 1
$ws.Cells.Item(28, 3).Value = 1
$ws.Cells.Item(28, 4).Value = 100
$ws.Cells.Item(28, 5).Value = 5.1
$ws.Cells.Item(28, 6).Value = 100
$ws.Cells.Item(29, 1).Value = "Oslo University Hospital"
$ws.Cells.Item(29, 2).Value = 102
$ws.Cells.Item(29, 3).Value = 50
$ws.Cells.Item(29, 4).Value = 49
$ws.Cells.Item(29, 5).Value = 39.5
$ws.Cells.Item(29, 6).Value = 58.59999999999999
$ws.Cells.Item(30, 1).Value = "Oulu University Hospital"
$ws.Cells.Item(30, 2).Value = 10
$ws.Cells.Item(30, 3).Value = 5
$ws.Cells.Item(30, 4).Value = 50
$ws.Cells.Item(30, 5).Value = 23.7
$ws.Cells.Item(30, 6).Value = 76.3
$ws.Cells.Item(31, 1).Value = "Sahlgrenska University Hospital"
$ws.Cells.Item(31, 2).Value = 42
$ws.Cells.Item(31, 3).Value = 17
$ws.Cells.Item(31, 4).Value = 40.5
$ws.Cells.Item(31, 5).Value = 27
$ws.Cells.Item(31, 6).Value = 55.50000000000001
$ws.Cells.Item(32, 1).Value = "Skane University Hospital"
$ws.Cells.Item(32, 2).Value = 23
$ws.Cells.Item(32, 3).Value = 14
$ws.Cells.Item(32, 4).Value = 60.9
$ws.Cells.Item(32, 5).Value = 40.8
$ws.Cells.Item(32, 6).Value = 77.8
$ws.Cells.Item(33, 1).Value = "St. Olav’s University Hospital"
$ws.Cells.Item(33, 2).Value = 24
$ws.Cells.Item(33, 3).Value = 8
$ws.Cells.Item(33, 4).Value = 33.3
$ws.Cells.Item(33, 5).Value = 18
$ws.Cells.Item(33, 6).Value = 53.3
$ws.Cells.Item(34, 1).Value = "Steno Diabetes Center Copenhagen"
$ws.Cells.Item(34, 2).Value = 13
$ws.Cells.Item(34, 3).Value = 9
$ws.Cells.Item(34, 4).Value = 69.2
$ws.Cells.Item(34, 5).Value = 42.4
$ws.Cells.Item(34, 6).Value = 87.3
$ws.Cells.Item(35, 1).Value = "Stockholm South General Hospital"
$ws.Cells.Item(35, 2).Value = 3
$ws.Cells.Item(35, 4).Value = 66.7
$ws.Cells.Item(35, 5).Value = 11.8
$ws.Cells.Item(35, 6).Value = 98.3
$ws.Cells.Item(36, 1).Value = "Tampere University Hospital"
$ws.Cells.Item(36, 2).Value = 24
$ws.Cells.Item(36, 3).Value = 16
$ws.Cells.Item(36, 4).Value = 66.7
$ws.Cells.Item(36, 5).Value = 46.7
$ws.Cells.Item(36, 6).Value = 82
$ws.Cells.Item(37, 1).Value = "The National University Hospital of Iceland"
$ws.Cells.Item(37, 2).Value = 5
$ws.Cells.Item(37, 3).Value = 2
$ws.Cells.Item(37, 4).Value = 40
$ws.Cells.Item(37, 5).Value = 7.1
$ws.Cells.Item(37, 6).Value = 76.9
$ws.Cells.Item(38, 1).Value = "Turku University Hospital"
$ws.Cells.Item(38, 2).Value = 50
$ws.Cells.Item(38, 4).Value = 40
$ws.Cells.Item(38, 5).Value = 27.6
$ws.Cells.Item(38, 6).Value = 53.8
$ws.Cells.Item(39, 1).Value = "UiT The Arctic University of Norway"
$ws.Cells.Item(39, 2).Value = 14
$ws.Cells.Item(39, 3).Value = 7
$ws.Cells.Item(39, 4).Value = 50
$ws.Cells.Item(39, 5).Value = 26.8
$ws.Cells.Item(39, 6).Value = 73.2
$ws.Cells.Item(40, 1).Value = "Umeå University"
$ws.Cells.Item(40, 2).Value = 42
$ws.Cells.Item(40, 3).Value = 20
$ws.Cells.Item(40, 4).Value = 47.6
$ws.Cells.Item(40, 5).Value = 33.4
$ws.Cells.Item(40, 6).Value = 62.3
$ws.Cells.Item(41, 1).Value = "University Hospital of North Norway"
$ws.Cells.Item(41, 2).Value = 17
$ws.Cells.Item(41, 3).Value = 6
$ws.Cells.Item(41, 4).Value = 35.3
$ws.Cells.Item(41, 5).Value = 17.3
$ws.Cells.Item(41, 6).Value = 58.7
$ws.Cells.Item(42, 1).Value = "University Hospital of Umeå"
$ws.Cells.Item(42, 2).Value = 2
$ws.Cells.Item(42, 3).Value = 1
$ws.Cells.Item(42, 4).Value = 50
$ws.Cells.Item(42, 5).Value = 2.6
$ws.Cells.Item(42, 6).Value = 97.39999999999999
$ws.Cells.Item(43, 1).Value = "University of Bergen"
$ws.Cells.Item(43, 2).Value = 31
$ws.Cells.Item(43, 3).Value = 16
$ws.Cells.Item(43, 4).Value = 51.6
$ws.Cells.Item(43, 5).Value = 34.8
$ws.Cells.Item(43, 6).Value = 68
$ws.Cells.Item(44, 1).Value = "University of Copenhagen"
$ws.Cells.Item(44, 2).Value = 99
$ws.Cells.Item(44, 3).Value = 43
$ws.Cells.Item(44, 4).Value = 43.4
$ws.Cells.Item(44, 5).Value = 34.1
$ws.Cells.Item(44, 6).Value = 53.3
$ws.Cells.Item(45, 1).Value = "University of Eastern Finland"
$ws.Cells.Item(45, 2).Value = 12
$ws.Cells.Item(45, 4).Value = 16.7
$ws.Cells.Item(45, 5).Value = 3
$ws.Cells.Item(45, 6).Value = 44.8
$ws.Cells.Item(46, 1).Value = "University of Helsinki"
$ws.Cells.Item(46, 2).Value = 21
$ws.Cells.Item(46, 3).Value = 11
$ws.Cells.Item(46, 4).Value = 52.4
$ws.Cells.Item(46, 5).Value = 32.4
$ws.Cells.Item(46, 6).Value = 71.7
$ws.Cells.Item(47, 1).Value = "University of Iceland"
$ws.Cells.Item(47, 2).Value = 5
$ws.Cells.Item(47, 3).Value = 2
$ws.Cells.Item(47, 5).Value = 7.1
$ws.Cells.Item(47, 6).Value = 76.9
$ws.Cells.Item(48, 1).Value = "University of Oslo"
$ws.Cells.Item(48, 2).Value = 23
$ws.Cells.Item(48, 3).Value = 15
$ws.Cells.Item(48, 4).Value = 65.2
$ws.Cells.Item(48, 5).Value = 44.9
$ws.Cells.Item(48, 6).Value = 81.2
$ws.Cells.Item(49, 1).Value = "University of Oulu"
$ws.Cells.Item(49, 2).Value = 25
$ws.Cells.Item(49, 3).Value = 10
$ws.Cells.Item(49, 4).Value = 40
$ws.Cells.Item(49, 5).Value = 23.4
$ws.Cells.Item(49, 6).Value = 59.3
$ws.Cells.Item(50, 1).Value = "University of Southern Denmark"
$ws.Cells.Item(50, 2).Value = 42
$ws.Cells.Item(50, 3).Value = 27
$ws.Cells.Item(50, 4).Value = 64.3
$ws.Cells.Item(50, 5).Value = 49.2
$ws.Cells.Item(50, 6).Value = 77
$ws.Cells.Item(51, 1).Value = "University of Tampere"
$ws.Cells.Item(51, 3).Value = 5
$ws.Cells.Item(51, 4).Value = 55.6
$ws.Cells.Item(51, 5).Value = 26.7
$ws.Cells.Item(51, 6).Value = 81.10000000000001
$ws.Cells.Item(52, 1).Value = "University of Turku"
$ws.Cells.Item(52, 2).Value = 20
$ws.Cells.Item(52, 3).Value = 10
$ws.Cells.Item(52, 4).Value = 50
$ws.Cells.Item(52, 5).Value = 29.9
$ws.Cells.Item(52, 6).Value = 70.1
$ws.Cells.Item(53, 1).Value = "Uppsala Academic Hospital"
$ws.Cells.Item(53, 2).Value = 9
$ws.Cells.Item(53, 3).Value = 3
$ws.Cells.Item(53, 4).Value = 33.3
$ws.Cells.Item(53, 5).Value = 12.1
$ws.Cells.Item(53, 6).Value = 64.60000000000001
$ws.Cells.Item(54, 1).Value = "Uppsala University"
$ws.Cells.Item(54, 2).Value = 51
$ws.Cells.Item(54, 3).Value = 22
$ws.Cells.Item(54, 4).Value = 43.1
$ws.Cells.Item(54, 5).Value = 30.5
$ws.Cells.Item(54, 6).Value = 56.7
$ws.Cells.Item(55, 1).Value = "Zealand University Hospital"
$ws.Cells.Item(55, 2).Value = 28
$ws.Cells.Item(55, 3).Value = 16
$ws.Cells.Item(55, 4).Value = 57.1
$ws.Cells.Item(55, 5).Value = 39.1
$ws.Cells.Item(55, 6).Value = 73.5
$ws.Cells.Item(56, 2).Value = 2112
$ws.Cells.Item(56, 5).Value = 49.6
